$wb = $excel.ActiveWorkbook

# --- Sheet2 ("Sheet2") -------------------------------------------------
# Loses tab-selection (Sheet1 takes over below), selection becomes E2:J2.
# Row2 -> Colombia lookup, Row3 -> Russia lookup, both FAIL.
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws2.Range("E2").Value = "codes=co"
$ws2.Range("J2").Value = "name=Colombia"
$ws2.Range("O2").Value = "FAIL"

$ws2.Range("E3").Value = "codes=rus"
$ws2.Range("J3").Value = "name=Russia"
$ws2.Range("O3").Value = "FAIL"

[void]$ws2.Activate()
[void]$ws2.Range("E2:J2").Select()

# --- Sheet3 ("Sheet3") -------------------------------------------------
# Selection becomes a single cell F2. Row2 -> India lookup, FAIL.
# Row 3 (the old "T2" test case) is removed.
$ws3 = $wb.Worksheets.Item("Sheet3")

$ws3.Range("E2").Value = "codes=in"
$ws3.Range("J2").Value = "name=India"
$ws3.Range("O2").Value = "FAIL"

[void]$ws3.Range("A3:O3").ClearContents()

[void]$ws3.Activate()
[void]$ws3.Range("F2").Select()

# --- Sheet1 ("Sheet1") -----------------------------------------------
# Becomes the active / tab-selected sheet, with a new selection of E2:J2.
# The REST test data changes from the generic "text=test" assertion to a
# country-code lookup for Norway, and the expected STATUS flips from
# PASS to FAIL. Row 3 (the old "T2" test case) is removed.
# Activated/selected LAST so it ends up as the workbook's active sheet.
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("E2").Value = "codes=no"
$ws1.Range("J2").Value = "name=Norway"
$ws1.Range("O2").Value = "FAIL"
$ws1.Range("O2").Style = "Normal"

[void]$ws1.Range("A3:O3").ClearContents()

[void]$ws1.Activate()
[void]$ws1.Range("E2:J2").Select()

Write-Output "edit complete"
